$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Despesa"
$ws.Range("B4").Value = "SERVIÇOS"
$ws.Range("C4").Value = 50
$ws.Range("D4").Value = "26/01/2025"

$ws.Range("A5").Value = "Receita"
$ws.Range("B5").Value = "ALUGUEL"
$ws.Range("C5").Value = 600
$ws.Range("D5").Value = "26/01/2025"
